# Update the "dSF" column (F) values on the active worksheet.
# These correspond to recalculated / repulled data values referenced
# in the commit message ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 1
    4  = -1
    5  = -6
    7  = -3
    8  = 2
    10 = 3
    11 = -1
    12 = 3
    13 = 2
    14 = 2
    15 = -4
    16 = -4
    17 = 1
    19 = -4
    21 = -3
    22 = -1
    23 = -2
    24 = -1
    25 = -5
    26 = 1
    28 = 4
    29 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
